$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renew the GDP mapping template: the "Index" unit variants are renamed to
# lowercase "index" variants (D4, D7, D10, D13, D16, D19).
$ws.Range("D4").Value  = "index"
$ws.Range("D7").Value  = "index, % YoY"
$ws.Range("D10").Value = "index, SA"
$ws.Range("D13").Value = "index, % MoM"
$ws.Range("D16").Value = "index, % YoY, SA"
$ws.Range("D19").Value = "index, % MoM, SA"

# Move the active selection from E22 to E20.
$ws.Range("E20").Select()
